# Emotion and hand gesture window added
# problems , Emotion links are incorrect, Single player not playlist player
#
# Fix the emotion -> song file links: the old links pointed at .m3u
# playlist files under a .\songs\ folder; the player now only needs a
# single .mp3 file per emotion, so replace each link cell with the
# matching mp3 filename.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 holds the link for each emotion in row 1 (A=angry, B=happy,
# C=sad, D=neutral). Update the song/link cells to the new single-file
# mp3 names.
$ws.Range("A2").Value = "angry.mp3"
$ws.Range("B2").Value = "happy.mp3"
$ws.Range("C2").Value = "sad.mp3"
$ws.Range("D2").Value = "neutral.mp3"

# Tweak column widths that were adjusted alongside the content fix.
$ws.Columns.Item(2).ColumnWidth = 21.86
$ws.Columns.Item(4).ColumnWidth = 25.71

# Scroll the view over a column and reselect the last-edited cell.
$ws.Range("D2").Select() | Out-Null
$excel.ActiveWindow.ScrollColumn = 2
